$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.431.74'
$ws.Range('E2').Value = '  -3.23%  '
$ws.Range('D3').Value = '2.470.19'
$ws.Range('E3').Value = '  -2.41%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '312.22'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = '95.04'
$ws.Range('E6').Value = '  -6.18%  '
$ws.Range('E7').Value = '  -2.64%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -4.42%  '
$ws.Range('D10').Value = '33.66'
$ws.Range('E10').Value = '  -6.06%  '
$ws.Range('E11').Value = '  -3.00%  '
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('D13').Value = '7.03'
$ws.Range('E13').Value = '  -4.46%  '
$ws.Range('D14').Value = '2.849.18'
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').Value = '2.436.42'
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('D16').Value = '14.90'
$ws.Range('E16').Value = '  -3.00%  '
$ws.Range('E17').Value = '  -3.57%  '
$ws.Range('D18').Value = '41.400.83'
$ws.Range('E18').Value = '  -3.20%  '
$ws.Range('D19').Value = '6.34'
$ws.Range('E19').Value = '  -4.99%  '
$ws.Range('D20').Value = '0.0₃0923'
$ws.Range('E20').Value = '  -3.22%  '
$ws.Range('D21').Value = '11.31'
$ws.Range('E21').Value = '  -9.19%  '
$ws.Range('D22').Value = '68.68'
$ws.Range('D23').Value = '237.47'
$ws.Range('E23').Value = '  -2.63%  '
$ws.Range('D24').Value = '2.75'
$ws.Range('E24').Value = '  -4.37%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').Value = '1.91'
$ws.Range('E26').Value = '  -6.46%  '
$ws.Range('D27').Value = '24.21'
$ws.Range('E27').Value = '  -5.25%  '
$ws.Range('E28').Value = '  -5.38%  '
$ws.Range('D29').Value = '9.65'
$ws.Range('E29').Value = '  -5.45%  '
$ws.Range('D30').Value = '36.81'
$ws.Range('E30').Value = '  -5.13%  '
$ws.Range('D31').Value = '151.90'
$ws.Range('E31').Value = '  -5.87%  '
$ws.Range('E32').Value = '  -5.41%  '
$ws.Range('D33').Value = '2.65'
$ws.Range('E33').Value = '  -4.82%  '
$ws.Range('D34').Value = '2.60'
$ws.Range('E34').Value = '  -2.64%  '
$ws.Range('E35').Value = '  -5.52%  '
$ws.Range('D36').Value = '3.06'
$ws.Range('E36').Value = '  -1.28%  '
$ws.Range('D37').Value = '1.89'
$ws.Range('E37').Value = '  -3.77%  '
$ws.Range('D38').Value = '17.03'
$ws.Range('E38').Value = '  -7.89%  '
$ws.Range('E39').Value = '  -2.71%  '
$ws.Range('D40').Value = '4.31'
$ws.Range('E40').Value = '  +3.25%  '
$ws.Range('E41').Value = '  -7.99%  '
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('D43').Value = '20.00'
$ws.Range('E43').Value = '  -9.54%  '
$ws.Range('D44').Value = '1.992.57'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('E45').Value = '  -4.35%  '
$ws.Range('D46').Value = '3.05'
$ws.Range('E46').Value = '  -8.37%  '
$ws.Range('D47').Value = '8.79'
$ws.Range('E47').Value = '  -2.66%  '
$ws.Range('D48').Value = '2.712.96'
$ws.Range('E48').Value = '  -2.06%  '
$ws.Range('D49').Value = '70.12'
$ws.Range('E49').Value = '  -3.37%  '
$ws.Range('D50').Value = '96.94'
$ws.Range('E50').Value = '  -4.22%  '
$ws.Range('D51').Value = '75.11'
$ws.Range('E51').Value = '  -5.77%  '
